$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Direccion" column (G1) is replaced by "Centro de formacion"
$ws.Range("G1").Value = "Centro de formacion"

# "Segundo Cargo" column (H1) is wiped out, leaving an empty, styled cell
$ws.Range("H1").Value = ""

# The last column (I1, formerly "Centro de formacion") is removed entirely,
# shrinking the sheet's used range/dimension down to column H
$ws.Range("I1").Clear()

# Restore the selection to G1 (was G4)
[void]$ws.Range("G1").Select()

# Add page setup info (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
